# Handle translations for block_inheritance in repository xlsx file.
#
# The repository sample sheet (Sheet1) gains two new example values for the
# "block_inheritance" column (Q) -- "nein" on the EinAmt example row (7) and
# "ja" on the GROUP_A example row (6) -- and the archival_value example on
# row 8 (column K) is translated from German "Archivwürdig" to English
# "archival worthy", which is also highlighted with a red box so it stands
# out as a to-be-translated / reference value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New example value in the "EinAmt" sample row -> block_inheritance = nein
$ws.Range("Q7").Value = "nein"

# Translate the archival_value example from German to English.
$ws.Range("K8").Value = "archival worthy"

# New example value in the "GROUP_A" sample row -> block_inheritance = ja
$ws.Range("Q6").Value = "ja"

# Highlight the translated archival_value example cell with a thin red box
# (white fill + thin red border) so translators can spot it quickly.
$ws.Range("K8").Borders.Color = 255
$ws.Range("K8").Interior.Color = 16777215

# Keep the header-row selection/scroll roughly where the author left it.
$ws.Activate()
$ws.Range("Q8").Select()

Write-Host "done"
